## Auto-generated script applying the diff to genx_signals.xlsx
## Checkpoint before follow-up message

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Copy-Format($src, $dst) {
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

function Set-TextValue($cell, $text) {
    # Force a literal text value even if it looks numeric/percent so Excel
    # does not silently convert it to a number with a new number format.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ==================== Sheet: Active Signals ====================
$ws1 = $wb.Worksheets.Item("Active Signals")

# Stash style donors (far away, same columns so the used dimension is not
# affected once we clear them again at the end).
$genericDonor = $ws1.Cells.Item(2, 1)          # plain data style (s=2)
$sellDonor    = $ws1.Cells.Item(200, 3)
$buyDonor     = $ws1.Cells.Item(201, 3)
$confDonor    = $ws1.Cells.Item(202, 8)
Copy-Format $ws1.Cells.Item(2, 3) $sellDonor   # original SELL style (s=3)
Copy-Format $ws1.Cells.Item(5, 3) $buyDonor    # original BUY style (s=4)
Copy-Format $ws1.Cells.Item(2, 8) $confDonor   # confidence text style (s=2)

# --- Row 2 ---
$ws1.Cells.Item(2, 1).Value = "2025-07-28 20:02"
$ws1.Cells.Item(2, 2).Value = "USDCHF"
$ws1.Cells.Item(2, 3).Value = "SELL"
Copy-Format $sellDonor $ws1.Cells.Item(2, 3)
$ws1.Cells.Item(2, 4).Value = 0.88091
$ws1.Cells.Item(2, 5).Value = 0.88518
$ws1.Cells.Item(2, 6).Value = 0.87119
$ws1.Cells.Item(2, 7).Value = 0.08
Set-TextValue $ws1.Cells.Item(2, 8) "87.0%"
Copy-Format $confDonor $ws1.Cells.Item(2, 8)
$ws1.Cells.Item(2, 9).Value = 2.27
$ws1.Cells.Item(2, 10).Value = "Active"

# --- Row 3 ---
$ws1.Cells.Item(3, 1).Value = "2025-07-28 20:04"
$ws1.Cells.Item(3, 2).Value = "NZDUSD"
$ws1.Cells.Item(3, 3).Value = "BUY"
Copy-Format $buyDonor $ws1.Cells.Item(3, 3)
$ws1.Cells.Item(3, 4).Value = 0.59179
$ws1.Cells.Item(3, 5).Value = 0.58828
$ws1.Cells.Item(3, 6).Value = 0.5992
$ws1.Cells.Item(3, 7).Value = 0.04
Set-TextValue $ws1.Cells.Item(3, 8) "84.0%"
Copy-Format $confDonor $ws1.Cells.Item(3, 8)
$ws1.Cells.Item(3, 9).Value = 2.11
$ws1.Cells.Item(3, 10).Value = "Active"

# --- Row 4 ---
$ws1.Cells.Item(4, 1).Value = "2025-07-28 19:53"
$ws1.Cells.Item(4, 2).Value = "USDJPY"
$ws1.Cells.Item(4, 3).Value = "BUY"
Copy-Format $buyDonor $ws1.Cells.Item(4, 3)
$ws1.Cells.Item(4, 4).Value = 150.15321
$ws1.Cells.Item(4, 5).Value = 149.91022
$ws1.Cells.Item(4, 6).Value = 150.6297
$ws1.Cells.Item(4, 7).Value = 0.04
Set-TextValue $ws1.Cells.Item(4, 8) "82.0%"
Copy-Format $confDonor $ws1.Cells.Item(4, 8)
$ws1.Cells.Item(4, 9).Value = 1.96
$ws1.Cells.Item(4, 10).Value = "Active"

# --- Row 5 ---
$ws1.Cells.Item(5, 1).Value = "2025-07-28 19:28"
$ws1.Cells.Item(5, 2).Value = "USDJPY"
$ws1.Cells.Item(5, 3).Value = "SELL"
Copy-Format $sellDonor $ws1.Cells.Item(5, 3)
$ws1.Cells.Item(5, 4).Value = 148.8192
$ws1.Cells.Item(5, 5).Value = 149.03316
$ws1.Cells.Item(5, 6).Value = 148.41725
$ws1.Cells.Item(5, 7).Value = 0.06
Set-TextValue $ws1.Cells.Item(5, 8) "83.0%"
Copy-Format $confDonor $ws1.Cells.Item(5, 8)
$ws1.Cells.Item(5, 9).Value = 1.88
$ws1.Cells.Item(5, 10).Value = "Active"

# --- Row 6 ---
$ws1.Cells.Item(6, 1).Value = "2025-07-28 20:02"
$ws1.Cells.Item(6, 2).Value = "USDCAD"
$ws1.Cells.Item(6, 3).Value = "SELL"
Copy-Format $sellDonor $ws1.Cells.Item(6, 3)
$ws1.Cells.Item(6, 4).Value = 1.37045
$ws1.Cells.Item(6, 5).Value = 1.37409
$ws1.Cells.Item(6, 6).Value = 1.36562
$ws1.Cells.Item(6, 7).Value = 0.01
Set-TextValue $ws1.Cells.Item(6, 8) "68.0%"
Copy-Format $confDonor $ws1.Cells.Item(6, 8)
$ws1.Cells.Item(6, 9).Value = 1.33
$ws1.Cells.Item(6, 10).Value = "Active"

# --- Row 7 ---
Copy-Format $genericDonor $ws1.Cells.Item(7, 1)
Copy-Format $genericDonor $ws1.Cells.Item(7, 2)
Copy-Format $genericDonor $ws1.Cells.Item(7, 4)
Copy-Format $genericDonor $ws1.Cells.Item(7, 5)
Copy-Format $genericDonor $ws1.Cells.Item(7, 6)
Copy-Format $genericDonor $ws1.Cells.Item(7, 7)
Copy-Format $genericDonor $ws1.Cells.Item(7, 9)
Copy-Format $genericDonor $ws1.Cells.Item(7, 10)
$ws1.Cells.Item(7, 1).Value = "2025-07-28 19:17"
$ws1.Cells.Item(7, 2).Value = "GBPUSD"
$ws1.Cells.Item(7, 3).Value = "BUY"
Copy-Format $buyDonor $ws1.Cells.Item(7, 3)
$ws1.Cells.Item(7, 4).Value = 1.27183
$ws1.Cells.Item(7, 5).Value = 1.2696
$ws1.Cells.Item(7, 6).Value = 1.28111
$ws1.Cells.Item(7, 7).Value = 0.07
Set-TextValue $ws1.Cells.Item(7, 8) "80.0%"
Copy-Format $confDonor $ws1.Cells.Item(7, 8)
$ws1.Cells.Item(7, 9).Value = 4.15
$ws1.Cells.Item(7, 10).Value = "Active"

# --- Row 8 ---
Copy-Format $genericDonor $ws1.Cells.Item(8, 1)
Copy-Format $genericDonor $ws1.Cells.Item(8, 2)
Copy-Format $genericDonor $ws1.Cells.Item(8, 4)
Copy-Format $genericDonor $ws1.Cells.Item(8, 5)
Copy-Format $genericDonor $ws1.Cells.Item(8, 6)
Copy-Format $genericDonor $ws1.Cells.Item(8, 7)
Copy-Format $genericDonor $ws1.Cells.Item(8, 9)
Copy-Format $genericDonor $ws1.Cells.Item(8, 10)
$ws1.Cells.Item(8, 1).Value = "2025-07-28 19:36"
$ws1.Cells.Item(8, 2).Value = "USDCHF"
$ws1.Cells.Item(8, 3).Value = "BUY"
Copy-Format $buyDonor $ws1.Cells.Item(8, 3)
$ws1.Cells.Item(8, 4).Value = 0.8793
$ws1.Cells.Item(8, 5).Value = 0.87649
$ws1.Cells.Item(8, 6).Value = 0.88683
$ws1.Cells.Item(8, 7).Value = 0.05
Set-TextValue $ws1.Cells.Item(8, 8) "77.0%"
Copy-Format $confDonor $ws1.Cells.Item(8, 8)
$ws1.Cells.Item(8, 9).Value = 2.68
$ws1.Cells.Item(8, 10).Value = "Active"

# remove the stashed donor cells again so the sheet dimension is unaffected
$sellDonor.Clear()
$buyDonor.Clear()
$confDonor.Clear()

# ==================== Sheet: Summary Dashboard ====================
$ws2 = $wb.Worksheets.Item("Summary Dashboard")
$ws2.Cells.Item(3, 2).Value = 15
$ws2.Cells.Item(4, 2).Value = 7
$ws2.Cells.Item(5, 2).Value = 6
$ws2.Cells.Item(6, 2).Value = 9
Set-TextValue $ws2.Cells.Item(7, 2) "78.3%"
Copy-Format $ws2.Cells.Item(3, 2) $ws2.Cells.Item(7, 2)
Set-TextValue $ws2.Cells.Item(8, 2) "2.54"
Copy-Format $ws2.Cells.Item(3, 2) $ws2.Cells.Item(8, 2)
Set-TextValue $ws2.Cells.Item(9, 2) "2025-07-28 19:45:30"
Copy-Format $ws2.Cells.Item(3, 2) $ws2.Cells.Item(9, 2)

# ==================== Sheet: Signal History ====================
$ws3 = $wb.Worksheets.Item("Signal History")
# --- Row 2 ---
$ws3.Cells.Item(2, 1).Value = "2025-07-28 20:02"
$ws3.Cells.Item(2, 2).Value = "USDCHF"
$ws3.Cells.Item(2, 3).Value = "BUY"
$ws3.Cells.Item(2, 4).Value = 0.88091
$ws3.Cells.Item(2, 5).Value = 0.88518
$ws3.Cells.Item(2, 6).Value = 0.87119
$ws3.Cells.Item(2, 7).Value = 0.08
$ws3.Cells.Item(2, 8).Value = 0.87
$ws3.Cells.Item(2, 9).Value = 2.27
$ws3.Cells.Item(2, 10).Value = "Active"

# --- Row 3 ---
$ws3.Cells.Item(3, 1).Value = "2025-07-28 20:04"
$ws3.Cells.Item(3, 2).Value = "NZDUSD"
$ws3.Cells.Item(3, 3).Value = "BUY"
$ws3.Cells.Item(3, 4).Value = 0.59179
$ws3.Cells.Item(3, 5).Value = 0.58828
$ws3.Cells.Item(3, 6).Value = 0.5992
$ws3.Cells.Item(3, 7).Value = 0.04
$ws3.Cells.Item(3, 8).Value = 0.84
$ws3.Cells.Item(3, 9).Value = 2.11
$ws3.Cells.Item(3, 10).Value = "Active"

# --- Row 4 ---
$ws3.Cells.Item(4, 1).Value = "2025-07-28 19:47"
$ws3.Cells.Item(4, 2).Value = "USDCAD"
$ws3.Cells.Item(4, 3).Value = "BUY"
$ws3.Cells.Item(4, 4).Value = 1.36364
$ws3.Cells.Item(4, 5).Value = 1.3615
$ws3.Cells.Item(4, 6).Value = 1.37179
$ws3.Cells.Item(4, 7).Value = 0.08
$ws3.Cells.Item(4, 8).Value = 0.85
$ws3.Cells.Item(4, 9).Value = 3.81
$ws3.Cells.Item(4, 10).Value = "Pending"

# --- Row 5 ---
$ws3.Cells.Item(5, 1).Value = "2025-07-28 19:55"
$ws3.Cells.Item(5, 2).Value = "USDCAD"
$ws3.Cells.Item(5, 3).Value = "SELL"
$ws3.Cells.Item(5, 4).Value = 1.36515
$ws3.Cells.Item(5, 5).Value = 1.36797
$ws3.Cells.Item(5, 6).Value = 1.36039
$ws3.Cells.Item(5, 7).Value = 0.09
$ws3.Cells.Item(5, 8).Value = 0.7
$ws3.Cells.Item(5, 9).Value = 1.69
$ws3.Cells.Item(5, 10).Value = "Pending"

# --- Row 6 ---
$ws3.Cells.Item(6, 1).Value = "2025-07-28 19:53"
$ws3.Cells.Item(6, 2).Value = "USDJPY"
$ws3.Cells.Item(6, 3).Value = "BUY"
$ws3.Cells.Item(6, 4).Value = 150.15321
$ws3.Cells.Item(6, 5).Value = 149.91022
$ws3.Cells.Item(6, 6).Value = 150.6297
$ws3.Cells.Item(6, 7).Value = 0.04
$ws3.Cells.Item(6, 8).Value = 0.82
$ws3.Cells.Item(6, 9).Value = 1.96
$ws3.Cells.Item(6, 10).Value = "Active"

# --- Row 7 ---
$ws3.Cells.Item(7, 1).Value = "2025-07-28 19:54"
$ws3.Cells.Item(7, 2).Value = "EURUSD"
$ws3.Cells.Item(7, 3).Value = "SELL"
$ws3.Cells.Item(7, 4).Value = 1.10416
$ws3.Cells.Item(7, 5).Value = 1.10896
$ws3.Cells.Item(7, 6).Value = 1.09589
$ws3.Cells.Item(7, 7).Value = 0.03
$ws3.Cells.Item(7, 8).Value = 0.76
$ws3.Cells.Item(7, 9).Value = 1.72
$ws3.Cells.Item(7, 10).Value = "Pending"

# --- Row 8 ---
$ws3.Cells.Item(8, 1).Value = "2025-07-28 19:47"
$ws3.Cells.Item(8, 2).Value = "GBPUSD"
$ws3.Cells.Item(8, 3).Value = "SELL"
$ws3.Cells.Item(8, 4).Value = 1.27152
$ws3.Cells.Item(8, 5).Value = 1.27443
$ws3.Cells.Item(8, 6).Value = 1.26464
$ws3.Cells.Item(8, 7).Value = 0.03
$ws3.Cells.Item(8, 8).Value = 0.68
$ws3.Cells.Item(8, 9).Value = 2.36
$ws3.Cells.Item(8, 10).Value = "Pending"

# --- Row 9 ---
$ws3.Cells.Item(9, 1).Value = "2025-07-28 19:56"
$ws3.Cells.Item(9, 2).Value = "USDCAD"
$ws3.Cells.Item(9, 3).Value = "SELL"
$ws3.Cells.Item(9, 4).Value = 1.36167
$ws3.Cells.Item(9, 5).Value = 1.36465
$ws3.Cells.Item(9, 6).Value = 1.35637
$ws3.Cells.Item(9, 7).Value = 0.03
$ws3.Cells.Item(9, 8).Value = 0.82
$ws3.Cells.Item(9, 9).Value = 1.78
$ws3.Cells.Item(9, 10).Value = "Pending"

# --- Row 10 ---
$ws3.Cells.Item(10, 1).Value = "2025-07-28 19:46"
$ws3.Cells.Item(10, 2).Value = "GBPUSD"
$ws3.Cells.Item(10, 3).Value = "BUY"
$ws3.Cells.Item(10, 4).Value = 1.27579
$ws3.Cells.Item(10, 5).Value = 1.27373
$ws3.Cells.Item(10, 6).Value = 1.28514
$ws3.Cells.Item(10, 7).Value = 0.09
$ws3.Cells.Item(10, 8).Value = 0.69
$ws3.Cells.Item(10, 9).Value = 4.55
$ws3.Cells.Item(10, 10).Value = "Filled"

# --- Row 11 ---
$ws3.Cells.Item(11, 1).Value = "2025-07-28 19:20"
$ws3.Cells.Item(11, 2).Value = "USDCHF"
$ws3.Cells.Item(11, 3).Value = "SELL"
$ws3.Cells.Item(11, 4).Value = 0.88015
$ws3.Cells.Item(11, 5).Value = 0.88296
$ws3.Cells.Item(11, 6).Value = 0.87033
$ws3.Cells.Item(11, 7).Value = 0.09
$ws3.Cells.Item(11, 8).Value = 0.93
$ws3.Cells.Item(11, 9).Value = 3.5
$ws3.Cells.Item(11, 10).Value = "Filled"

# --- Row 12 ---
$ws3.Cells.Item(12, 1).Value = "2025-07-28 19:28"
$ws3.Cells.Item(12, 2).Value = "USDJPY"
$ws3.Cells.Item(12, 3).Value = "SELL"
$ws3.Cells.Item(12, 4).Value = 148.8192
$ws3.Cells.Item(12, 5).Value = 149.03316
$ws3.Cells.Item(12, 6).Value = 148.41725
$ws3.Cells.Item(12, 7).Value = 0.06
$ws3.Cells.Item(12, 8).Value = 0.83
$ws3.Cells.Item(12, 9).Value = 1.88
$ws3.Cells.Item(12, 10).Value = "Active"

# --- Row 13 ---
$ws3.Cells.Item(13, 1).Value = "2025-07-28 20:02"
$ws3.Cells.Item(13, 2).Value = "USDCAD"
$ws3.Cells.Item(13, 3).Value = "SELL"
$ws3.Cells.Item(13, 4).Value = 1.37045
$ws3.Cells.Item(13, 5).Value = 1.37409
$ws3.Cells.Item(13, 6).Value = 1.36562
$ws3.Cells.Item(13, 7).Value = 0.01
$ws3.Cells.Item(13, 8).Value = 0.68
$ws3.Cells.Item(13, 9).Value = 1.33
$ws3.Cells.Item(13, 10).Value = "Active"

# --- Row 14 ---
$ws3.Cells.Item(14, 1).Value = "2025-07-28 19:53"
$ws3.Cells.Item(14, 2).Value = "AUDUSD"
$ws3.Cells.Item(14, 3).Value = "SELL"
$ws3.Cells.Item(14, 4).Value = 0.65505
$ws3.Cells.Item(14, 5).Value = 0.65814
$ws3.Cells.Item(14, 6).Value = 0.64799
$ws3.Cells.Item(14, 7).Value = 0.07
$ws3.Cells.Item(14, 8).Value = 0.7
$ws3.Cells.Item(14, 9).Value = 2.28
$ws3.Cells.Item(14, 10).Value = "Pending"

# --- Row 15 ---
$ws3.Cells.Item(15, 1).Value = "2025-07-28 19:17"
$ws3.Cells.Item(15, 2).Value = "GBPUSD"
$ws3.Cells.Item(15, 3).Value = "BUY"
$ws3.Cells.Item(15, 4).Value = 1.27183
$ws3.Cells.Item(15, 5).Value = 1.2696
$ws3.Cells.Item(15, 6).Value = 1.28111
$ws3.Cells.Item(15, 7).Value = 0.07
$ws3.Cells.Item(15, 8).Value = 0.8
$ws3.Cells.Item(15, 9).Value = 4.15
$ws3.Cells.Item(15, 10).Value = "Active"

# --- Row 16 ---
$ws3.Cells.Item(16, 1).Value = "2025-07-28 19:36"
$ws3.Cells.Item(16, 2).Value = "USDCHF"
$ws3.Cells.Item(16, 3).Value = "BUY"
$ws3.Cells.Item(16, 4).Value = 0.8793
$ws3.Cells.Item(16, 5).Value = 0.87649
$ws3.Cells.Item(16, 6).Value = 0.88683
$ws3.Cells.Item(16, 7).Value = 0.05
$ws3.Cells.Item(16, 8).Value = 0.77
$ws3.Cells.Item(16, 9).Value = 2.68
$ws3.Cells.Item(16, 10).Value = "Active"

